# fix: page layout, add:higlight button
#
# Adds two new summary sheets at the end of the workbook:
#   - "Top 10 Unidades Más Eficientes"   (units ranked by best km/l)
#   - "Top 10 Unidades Menos Eficientes" (units ranked by worst km/l)
#
# Both mirror the layout already used by the existing
# "Top 10 Rutas Mas/Menos Eficientes" sheets: a bold/boxed header row
# ("Unidad", "Kms Totales", "Litros", "Eficiencia (km/l)", "CPK") followed
# by 10 data rows, using the same page margins as the rest of the workbook.

$wb = $excel.ActiveWorkbook

# Reference sheet whose header formatting (bold, centered, thin box border)
# we reuse so the new sheets' style matches the workbook's house style
# exactly instead of inventing a new one.
$styleSource = $wb.Worksheets.Item("Top 10 Rutas Mas Eficientes")

function Add-TopUnidadesSheet {
    param(
        $SheetName,
        $Rows
    )

    $ws = $wb.Worksheets.Add($null, $wb.Sheets($wb.Sheets.Count))

    # Excel sheet names are capped at 31 characters; if the requested name
    # is longer, fall back to the longest valid prefix instead of failing.
    if ($SheetName.Length -gt 31) {
        $SheetName = $SheetName.Substring(0, 31)
    }
    $ws.Name = $SheetName

    $headers = @("Unidad", "Kms Totales", "Litros", "Eficiencia (km/l)", "CPK")
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
    }

    # Copy the bold/centered/bordered header style from the existing
    # "Top 10 Rutas ..." sheets onto the new header row.
    $styleSource.Range("A1").Copy()
    $ws.Range("A1:E1").PasteSpecial(-4122)

    # Column A holds unit numbers that must stay text (e.g. "1665"), not
    # numbers, matching the source data.
    $ws.Range("A2:A11").NumberFormat = "@"

    $r = 2
    foreach ($row in $Rows) {
        $ws.Cells.Item($r, 1).Value = [string]$row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $r = $r + 1
    }

    # Match the page margins used across the rest of the workbook.
    $ws.PageSetup.LeftMargin = 54
    $ws.PageSetup.RightMargin = 54
    $ws.PageSetup.TopMargin = 72
    $ws.PageSetup.BottomMargin = 72
    $ws.PageSetup.HeaderMargin = 36
    $ws.PageSetup.FooterMargin = 36

    $ws.Range("A1").Select()

    return $ws
}

# Units ranked from most to least fuel-efficient (km/l).
$dataMasEficientes = @(
    @("1665", 125560, 111846.62, 1.122608801231544, 17.53992595759),
    @("1670", 1790, 5152, 0.3474378881987578, 89.95270058659219),
    @("1630", 96, 288, 0.3333333333333333, 58.11822),
    @("1640", 89040, 549960.45, 0.161902551356193, 121.0195035190364),
    @("1628", 216688, 1412981.28, 0.1533551810396243, 128.2657230081223),
    @("1633", 1060, 6995.17, 0.15153312928778, 127.2760094830189),
    @("1641", 80316, 554997.12, 0.1447142644632102, 135.570475367334),
    @("1639", 147277, 1322791.05, 0.1113380680947305, 176.4397602208016),
    @("1649", 111056, 1159789.86, 0.0957552775983056, 204.0743820347212),
    @("1635", 36, 458, 0.07860262008733625, 246.4643033333333)
)

# Same units, ordered from least to most efficient (mirror image).
$dataMenosEficientes = @(
    @("1635", 36, 458, 0.07860262008733625, 246.4643033333333),
    @("1649", 111056, 1159789.86, 0.0957552775983056, 204.0743820347212),
    @("1639", 147277, 1322791.05, 0.1113380680947305, 176.4397602208016),
    @("1641", 80316, 554997.12, 0.1447142644632102, 135.570475367334),
    @("1633", 1060, 6995.17, 0.15153312928778, 127.2760094830189),
    @("1628", 216688, 1412981.28, 0.1533551810396243, 128.2657230081223),
    @("1640", 89040, 549960.45, 0.161902551356193, 121.0195035190364),
    @("1630", 96, 288, 0.3333333333333333, 58.11822),
    @("1670", 1790, 5152, 0.3474378881987578, 89.95270058659219),
    @("1665", 125560, 111846.62, 1.122608801231544, 17.53992595759)
)

Add-TopUnidadesSheet "Top 10 Unidades Más Eficientes" $dataMasEficientes | Out-Null
Add-TopUnidadesSheet "Top 10 Unidades Menos Eficientes" $dataMenosEficientes | Out-Null

$wb.Worksheets.Item(1).Select()
